$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 2440.818
$ws.Range("I51").Value = 1699
$ws.Range("J51").Value = 2719
$ws.Range("K51").Value = 1699
$ws.Range("L51").Value = 2719
$ws.Range("M51").Value = -1215
$ws.Range("N51").Value = -3687
# Row 112
$ws.Range("H112").Value = 20409520
$ws.Range("J112").Value = 24846230
$ws.Range("L112").Value = 74538690
$ws.Range("N112").Value = -74540906
# Row 125
$ws.Range("H125").Value = 2378.3845
$ws.Range("I125").Value = 3390.375
$ws.Range("J125").Value = 759.2
$ws.Range("K125").Value = 30513.375
$ws.Range("L125").Value = 6832.8
$ws.Range("M125").Value = -28053.375
$ws.Range("N125").Value = -11752.8
# Row 138
$ws.Range("H138").Value = 3120.674
$ws.Range("I138").Value = 1512.3429
$ws.Range("J138").Value = 4163.1113
$ws.Range("K138").Value = 4537.028700000001
$ws.Range("L138").Value = 12489.3339
$ws.Range("M138").Value = 602.9712999999992
$ws.Range("N138").Value = -22769.3339
# Row 141
$ws.Range("H141").Value = 3102.6333
$ws.Range("I141").Value = 2429.3684
$ws.Range("K141").Value = 7288.1052
$ws.Range("M141").Value = -2108.1052

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13319.798
$ws.Range("I32").Value = 10845.889
$ws.Range("J32").Value = 19916.889
$ws.Range("K32").Value = 10845.889
$ws.Range("L32").Value = 19916.889
$ws.Range("M32").Value = -10558.889
$ws.Range("N32").Value = -20490.889
# Row 101
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
# Row 139
$ws.Range("H139").Value = 65787.86
$ws.Range("J139").Value = 65787.86
$ws.Range("L139").Value = 65787.86
$ws.Range("N139").Value = -76067.86

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 59
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents() | Out-Null
# Row 81
$ws.Range("H81").Value = 41924.75
$ws.Range("J81").Value = 41924.75
$ws.Range("L81").Value = 41924.75
$ws.Range("N81").Value = -44046.75
# Row 84
$ws.Range("H84").Value = 41924.75
$ws.Range("J84").Value = 41924.75
$ws.Range("L84").Value = 125774.25
$ws.Range("N84").Value = -136382.25
# Row 132
$ws.Range("H132").Value = 54780
$ws.Range("J132").Value = 54780
$ws.Range("L132").Value = 54780
$ws.Range("N132").Value = -64900

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8625168
$ws.Range("I31").Value = 1469.2424
$ws.Range("J31").Value = 20008450
$ws.Range("K31").Value = 1469.2424
$ws.Range("L31").Value = 20008450
$ws.Range("M31").Value = -1174.2424
$ws.Range("N31").Value = -20009040
# Row 34
$ws.Range("H34").Value = 8625168
$ws.Range("I34").Value = 1469.2424
$ws.Range("J34").Value = 20008450
$ws.Range("K34").Value = 1469.2424
$ws.Range("L34").Value = 20008450
$ws.Range("M34").Value = -1267.2424
$ws.Range("N34").Value = -20008854
# Row 58
$ws.Range("H58").Value = 3760007.8
$ws.Range("I58").Value = 4167679.2
$ws.Range("J58").Value = 1430457.2
$ws.Range("K58").Value = 4167679.2
$ws.Range("L58").Value = 1430457.2
$ws.Range("M58").Value = -4167476.2
$ws.Range("N58").Value = -1430863.2
# Row 105
$ws.Range("H105").Value = 20835718
$ws.Range("I105").Value = 25643384
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 25643384
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -25641637
$ws.Range("N105").Value = -5994
# Row 136
$ws.Range("H136").Value = 3760007.8
$ws.Range("I136").Value = 4167679.2
$ws.Range("J136").Value = 1430457.2
$ws.Range("K136").Value = 12503037.6
$ws.Range("L136").Value = 4291371.6
$ws.Range("M136").Value = -12500487.6
$ws.Range("N136").Value = -4296471.6
# Row 140
$ws.Range("H140").Value = 28464.516
$ws.Range("J140").Value = 28464.516
$ws.Range("L140").Value = 28464.516
$ws.Range("N140").Value = -38824.516

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2274269.5
$ws.Range("I5").Value = 623.04
$ws.Range("K5").Value = 1869.12
$ws.Range("M5").Value = -1757.12
# Row 38
$ws.Range("H38").Value = 4166942.2
$ws.Range("I38").Value = 6250324.5
$ws.Range("K38").Value = 18750973.5
$ws.Range("M38").Value = -18750626.5
# Row 113
$ws.Range("H113").Value = 2069454.6
$ws.Range("I113").Value = 16666934
$ws.Range("J113").Value = 385130.06
$ws.Range("K113").Value = 50000802
$ws.Range("L113").Value = 1155390.18
$ws.Range("M113").Value = -49998632
$ws.Range("N113").Value = -1159730.18
# Row 131
$ws.Range("H131").Value = 2084309.8
$ws.Range("I131").Value = 5556195.5
$ws.Range("J131").Value = 1178.3334
$ws.Range("K131").Value = 16668586.5
$ws.Range("L131").Value = 3535.0002
$ws.Range("M131").Value = -16663546.5
$ws.Range("N131").Value = -13615.0002
# Row 132
$ws.Range("H132").Value = 4833344
$ws.Range("I132").Value = 1511
$ws.Range("J132").Value = 7939522.5
$ws.Range("K132").Value = 13599
$ws.Range("L132").Value = 71455702.5
$ws.Range("M132").Value = -11069
$ws.Range("N132").Value = -71460762.5
# Row 135
$ws.Range("H135").Value = 2274269.5
$ws.Range("I135").Value = 623.04
$ws.Range("K135").Value = 5607.36
$ws.Range("M135").Value = -3072.36

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 14688.125
$ws.Range("I80").Value = 14688.125
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 14688.125
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -13690.125
$ws.Range("N80").ClearContents() | Out-Null
# Row 83
$ws.Range("H83").Value = 14688.125
$ws.Range("I83").Value = 14688.125
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 73440.625
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -68448.625
$ws.Range("N83").ClearContents() | Out-Null
# Row 140
$ws.Range("H140").Value = 37701.633
$ws.Range("J140").Value = 37701.633
$ws.Range("L140").Value = 37701.633
$ws.Range("N140").Value = -48061.633

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 654594.1
$ws.Range("I82").Value = 1112677.4
$ws.Range("J82").Value = 139250.5
$ws.Range("K82").Value = 1112677.4
$ws.Range("L82").Value = 139250.5
$ws.Range("M82").Value = -1112316.4
$ws.Range("N82").Value = -139972.5
# Row 85
$ws.Range("H85").Value = 654594.1
$ws.Range("I85").Value = 1112677.4
$ws.Range("J85").Value = 139250.5
$ws.Range("K85").Value = 1112677.4
$ws.Range("L85").Value = 139250.5
$ws.Range("M85").Value = -1111429.4
$ws.Range("N85").Value = -141746.5
# Row 103
$ws.Range("H103").Value = 43000
$ws.Range("J103").Value = 43000
$ws.Range("L103").Value = 43000
$ws.Range("N103").Value = -45344
# Row 136
$ws.Range("H136").Value = 10004.087
$ws.Range("I136").Value = 4280.909
$ws.Range("J136").Value = 15250.333
$ws.Range("K136").Value = 12842.727
$ws.Range("L136").Value = 45750.999
$ws.Range("M136").Value = -10292.727
$ws.Range("N136").Value = -50850.999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 101
$ws.Range("H101").Value = 16200.667
$ws.Range("J101").Value = 16200.667
$ws.Range("L101").Value = 16200.667
$ws.Range("N101").Value = -22690.667
# Row 103
$ws.Range("H103").Value = 33375
$ws.Range("J103").Value = 33375
$ws.Range("L103").Value = 33375
$ws.Range("N103").Value = -35719
# Row 132
$ws.Range("H132").Value = 1520.15
$ws.Range("I132").Value = 687.55
$ws.Range("J132").Value = 2352.75
$ws.Range("K132").Value = 2062.65
$ws.Range("L132").Value = 7058.25
$ws.Range("M132").Value = 467.3500000000004
$ws.Range("N132").Value = -12118.25

Write-Host "Applied all Ixion Profits updates"